# Updated cryptos list (prices / 1h volume %) scraped on 2024-06-15.
# All written values are prefixed with a literal apostrophe so Excel
# treats numeric-looking strings (e.g. "0.137", "66.303.94") as text,
# matching the inline-string cells already used throughout the sheet.
# The Style is reset to Normal right after so the quote-prefix marker
# does not leave a stray cell style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.303.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.81%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.536.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.78%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.06%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''607.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.02%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''143.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.30%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''3.534.86'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.68%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.02%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  +0.53%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  +1.30%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.137'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -3.81%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.410'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -2.72%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''4.138.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = '''  -3.98%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''30.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -5.32%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''3.539.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.63%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''66.391.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -0.88%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  -0.67%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''10.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +2.11%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  -3.73%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''14.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.45%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''425.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -2.65%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.601'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -1.07%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''78.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -1.09%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''3.680.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.72%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  +0.03%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  +0.39%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''8.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -1.22%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -5.64%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -1.62%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -0.32%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''1.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -7.54%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -4.36%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -0.84%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''3.529.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.65%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -0.03%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''1.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -2.80%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = '''Aptos'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = '''7.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -2.64%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = '''NEARProtocol'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = '''5.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -5.18%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -0.09%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''172.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -1.46%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.0855'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -4.10%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -4.04%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.893'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -0.05%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -7.49%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''45.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.31%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -1.97%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''26.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -7.17%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -1.21%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -4.21%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.945'
$ws.Range("D51").Style = "Normal"
